# Python Lesson 3 Functions.pptx - insert new "scope" section-header slide
#
# The original slide at (1-based) position 22 ("LIbraries") stays put; a
# brand new Section-Header slide ("scope") is inserted immediately before
# it, so it becomes the new position 22 and every slide that used to be
# at position >= 22 shifts down by one.

$p = $ppt.ActivePresentation

# "Section Header" is CustomLayout index 3 on the (single) slide master -
# matches ppt/slideLayouts/slideLayout3.xml, the same layout used by the
# neighboring "LIbraries"/"Libraries" section-header slides.
$sectionHeaderLayout = $p.SlideMaster.CustomLayouts.Item(3)

# Insert the new slide at position 22 (existing slides 22..26 shift to 23..27).
$newSlide = $p.Slides.AddSlide(22, $sectionHeaderLayout)

# Title placeholder -> "scope"; body placeholder is left empty, matching
# the target markup (endParaRPr only, no run).
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "scope"
